# Extend the "pbsum" (Invoicing Periods Workload) table: Alojzy now has
# three monthly rows (Jan/Feb/Mar, each with its own Lower/Upper bound)
# and the other four experts (Eugeniusz, Ignacy, Pafnucy, Romuald) keep
# their single January row, shifted down to make room.

$wb = $excel.ActiveWorkbook

# --- "period" sheet: just move the selection, no data change ---
$wsPeriod = $wb.Worksheets.Item("period")
$wsPeriod.Range("A4").Select() | Out-Null

# --- "pbsum" sheet: becomes the active sheet/tab ---
$ws = $wb.Worksheets.Item("pbsum")
$ws.Activate() | Out-Null

# Insert two fresh rows at row 7 (pushes the trailing blank spacer rows
# from 7-9 down to 9-11, and auto-extends the AND(...) ranges/dimension).
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# Row 2 (Alojzy / 25.Jan / 0 / 200) is unchanged.

# Row 3: Alojzy / 25.Feb / 10 / 190
$ws.Cells.Item(3, 1).Value = "Alojzy"
$ws.Cells.Item(3, 2).Value = "25.Feb"
$ws.Cells.Item(3, 3).Value = 10
$ws.Cells.Item(3, 4).Value = 190

# Row 4: Alojzy / 25.Mar / 20 / 180
$ws.Cells.Item(4, 1).Value = "Alojzy"
$ws.Cells.Item(4, 2).Value = "25.Mar"
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 4).Value = 180

# Row 5: Eugeniusz / 25.Jan / 0 / 200
$ws.Cells.Item(5, 1).Value = "Eugeniusz"
$ws.Cells.Item(5, 2).Value = "25.Jan"
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 200

# Row 6: Ignacy / 25.Jan / 0 / 200
$ws.Cells.Item(6, 1).Value = "Ignacy"
$ws.Cells.Item(6, 2).Value = "25.Jan"
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 200

# Row 7 (new): Pafnucy / 25.Jan / 0 / 200
$ws.Cells.Item(7, 1).Value = "Pafnucy"
$ws.Cells.Item(7, 2).Value = "25.Jan"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 200
$ws.Cells.Item(7, 5).Formula = "=COUNTIF(expert!`$A`$2:`$A`$921, A7) > 0"
$ws.Cells.Item(7, 6).Formula = "=COUNTIF(period!`$A`$2:`$A`$1000, B7) > 0"

# Row 8 (new): Romuald / 25.Jan / 0 / 200
$ws.Cells.Item(8, 1).Value = "Romuald"
$ws.Cells.Item(8, 2).Value = "25.Jan"
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 200
$ws.Cells.Item(8, 5).Formula = "=COUNTIF(expert!`$A`$2:`$A`$921, A8) > 0"
$ws.Cells.Item(8, 6).Formula = "=COUNTIF(period!`$A`$2:`$A`$1000, B8) > 0"

# Final selection on the now-active "pbsum" sheet
$ws.Range("F6").Select() | Out-Null
